$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a "model" label above the model-name header row
$ws.Range("B4").Value = "model"

# Insert a new blank row before row 6 (before "FSM self acc"),
# shifting existing rows 6-15 down to rows 7-16
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new "X-115" label cells
$ws.Range("D6").Value = "X-115"
$ws.Range("E6").Value = "X-115"

# Update values in row 7 ("FSM self acc", previously row 6)
$ws.Range("D7").Value = 0.9182
$ws.Range("E7").ClearContents()

# Update values in row 8 ("FSM self L2", previously row 7)
$ws.Range("D8").Value = 0.2078
$ws.Range("E8").ClearContents()

# Update values in row 10 ("Transfer to resnet", previously row 9)
$ws.Range("C10").Value = 0.2091

# Update values in row 11 ("Transfer to vgg", previously row 10)
$ws.Range("C11").Value = 0.3091

# Insert another blank row before "attack from (-1,1) acc"
# (currently at row 12 after the first insertion), shifting rows 12-16 down to 13-17
$ws.Rows.Item(12).Insert()

# Add new E column values for the attack rows (now rows 13 and 14)
$ws.Range("E13").Value = 0.3818
$ws.Range("E14").Value = 25.37

# Update the selection to match the recorded cursor position
$ws.Range("E10").Select()
